# Generate Report for Handback
#
# This updates the localization-status workbook so that the zh-cn and
# de-de sheets reflect that the handback from the vendor is complete:
#   - the "Status" column changes from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - two new columns get populated for each file row:
#       F "Latest Target File"   - the (markdown) source file name
#       G "Latest Handback File" - the localized handback file name
#     both rendered as hyperlinks, matching the look of the existing
#     hyperlinked columns
#   - the "Latest Handback DateTime" column (H) is stamped with the
#     real handback timestamp instead of the "0001-01-01 00:00:00"
#     placeholder (each language got handed back at its own time)

$wb = $excel.ActiveWorkbook

# Color used by the workbook's "HyperLink" cell style (RGB 6495ED,
# cornflower blue), expressed as the BGR-encoded long that Excel's
# Font.Color expects.
$hyperlinkColor = 15570276

$srcMdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/8e98f2051904923de05efd902c4af3697e35d9b1/e2e/29e79f51-6ede-4853-a79d-4cea48aefdf7.md"
$srcMdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/8e98f2051904923de05efd902c4af3697e35d9b1/e2e/69bc8315-b512-49ab-a3b9-5d471a9f1a0e.md"

$statusHandedBack = "Handed back: in sync with en-US"

# -----------------------------------------------------------------------
# Sheet "zh-cn"
# -----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusHandedBack
$wsZh.Range("C3").Value = $statusHandedBack

$wsZh.Range("F2").Value = "29e79f51-6ede-4853-a79d-4cea48aefdf7.md"
$wsZh.Range("G2").Value = "29e79f51-6ede-4853-a79d-4cea48aefdf7.33aa57355ad5153a5270f755ce14a331139f09e6.zh-cn.xlf"
$wsZh.Range("F3").Value = "69bc8315-b512-49ab-a3b9-5d471a9f1a0e.md"
$wsZh.Range("G3").Value = "69bc8315-b512-49ab-a3b9-5d471a9f1a0e.3c07634513735dd3901ca696129c0c39357ebf3d.zh-cn.xlf"

$wsZh.Range("H2").Value = "2016-03-18 08:46:22"
$wsZh.Range("H3").Value = "2016-03-18 08:46:22"

# Rebuild every hyperlink on the sheet so the final order follows the
# natural reading order: A2,B2,D2,F2,G2,A3,B3,D3,F3,G3
$wsZh.Range("A1").Hyperlinks.Delete()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $srcMdUrl1, "", "", "29e79f51-6ede-4853-a79d-4cea48aefdf7.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $srcMdUrl1, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a066214b83afa2279d9a41d28bdd90e0dc0a912e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/29e79f51-6ede-4853-a79d-4cea48aefdf7.33aa57355ad5153a5270f755ce14a331139f09e6.zh-cn.xlf", "", "", "29e79f51-6ede-4853-a79d-4cea48aefdf7.33aa57355ad5153a5270f755ce14a331139f09e6.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $srcMdUrl1, "", "", "29e79f51-6ede-4853-a79d-4cea48aefdf7.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a066214b83afa2279d9a41d28bdd90e0dc0a912e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/29e79f51-6ede-4853-a79d-4cea48aefdf7.33aa57355ad5153a5270f755ce14a331139f09e6.zh-cn.xlf", "", "", "29e79f51-6ede-4853-a79d-4cea48aefdf7.33aa57355ad5153a5270f755ce14a331139f09e6.zh-cn.xlf")

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $srcMdUrl2, "", "", "69bc8315-b512-49ab-a3b9-5d471a9f1a0e.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $srcMdUrl2, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a066214b83afa2279d9a41d28bdd90e0dc0a912e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/69bc8315-b512-49ab-a3b9-5d471a9f1a0e.3c07634513735dd3901ca696129c0c39357ebf3d.zh-cn.xlf", "", "", "69bc8315-b512-49ab-a3b9-5d471a9f1a0e.3c07634513735dd3901ca696129c0c39357ebf3d.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $srcMdUrl2, "", "", "69bc8315-b512-49ab-a3b9-5d471a9f1a0e.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a066214b83afa2279d9a41d28bdd90e0dc0a912e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/69bc8315-b512-49ab-a3b9-5d471a9f1a0e.3c07634513735dd3901ca696129c0c39357ebf3d.zh-cn.xlf", "", "", "69bc8315-b512-49ab-a3b9-5d471a9f1a0e.3c07634513735dd3901ca696129c0c39357ebf3d.zh-cn.xlf")

# Re-apply the workbook's hyperlink look (underline + cornflower blue)
# to every linked cell on the sheet, since adding a hyperlink resets
# the cell to Excel's generic built-in "Hyperlink" style.
foreach ($addr in @("A2","B2","D2","F2","G2","A3","B3","D3","F3","G3")) {
    $c = $wsZh.Range($addr)
    $c.Font.Underline = 2
    $c.Font.Color = $hyperlinkColor
}

# -----------------------------------------------------------------------
# Sheet "de-de"
# -----------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusHandedBack
$wsDe.Range("C3").Value = $statusHandedBack

$wsDe.Range("F2").Value = "29e79f51-6ede-4853-a79d-4cea48aefdf7.md"
$wsDe.Range("G2").Value = "29e79f51-6ede-4853-a79d-4cea48aefdf7.33aa57355ad5153a5270f755ce14a331139f09e6.de-de.xlf"
$wsDe.Range("F3").Value = "69bc8315-b512-49ab-a3b9-5d471a9f1a0e.md"
$wsDe.Range("G3").Value = "69bc8315-b512-49ab-a3b9-5d471a9f1a0e.3c07634513735dd3901ca696129c0c39357ebf3d.de-de.xlf"

$wsDe.Range("H2").Value = "2016-03-18 08:46:30"
$wsDe.Range("H3").Value = "2016-03-18 08:46:30"

$wsDe.Range("A1").Hyperlinks.Delete()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $srcMdUrl1, "", "", "29e79f51-6ede-4853-a79d-4cea48aefdf7.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $srcMdUrl1, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8f2b64f41508255d80e8836c8296aa6cfd3860ac/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/29e79f51-6ede-4853-a79d-4cea48aefdf7.33aa57355ad5153a5270f755ce14a331139f09e6.de-de.xlf", "", "", "29e79f51-6ede-4853-a79d-4cea48aefdf7.33aa57355ad5153a5270f755ce14a331139f09e6.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $srcMdUrl1, "", "", "29e79f51-6ede-4853-a79d-4cea48aefdf7.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8f2b64f41508255d80e8836c8296aa6cfd3860ac/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/29e79f51-6ede-4853-a79d-4cea48aefdf7.33aa57355ad5153a5270f755ce14a331139f09e6.de-de.xlf", "", "", "29e79f51-6ede-4853-a79d-4cea48aefdf7.33aa57355ad5153a5270f755ce14a331139f09e6.de-de.xlf")

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $srcMdUrl2, "", "", "69bc8315-b512-49ab-a3b9-5d471a9f1a0e.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $srcMdUrl2, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8f2b64f41508255d80e8836c8296aa6cfd3860ac/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/69bc8315-b512-49ab-a3b9-5d471a9f1a0e.3c07634513735dd3901ca696129c0c39357ebf3d.de-de.xlf", "", "", "69bc8315-b512-49ab-a3b9-5d471a9f1a0e.3c07634513735dd3901ca696129c0c39357ebf3d.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $srcMdUrl2, "", "", "69bc8315-b512-49ab-a3b9-5d471a9f1a0e.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8f2b64f41508255d80e8836c8296aa6cfd3860ac/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/69bc8315-b512-49ab-a3b9-5d471a9f1a0e.3c07634513735dd3901ca696129c0c39357ebf3d.de-de.xlf", "", "", "69bc8315-b512-49ab-a3b9-5d471a9f1a0e.3c07634513735dd3901ca696129c0c39357ebf3d.de-de.xlf")

foreach ($addr in @("A2","B2","D2","F2","G2","A3","B3","D3","F3","G3")) {
    $c = $wsDe.Range($addr)
    $c.Font.Underline = 2
    $c.Font.Color = $hyperlinkColor
}

$wb.Save()
